$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All")
$lo = $ws.ListObjects.Item("Table18")

# Add the two new table columns (this also extends the table range and
# grows the sheet dimension automatically).
$col1 = $lo.ListColumns.Add()
$ws.Range("R1").Value = "tug_st_time"
$col2 = $lo.ListColumns.Add()
$ws.Range("S1").Value = "tug_dt_time"

# Populate the new TUG single-task / dual-task time columns for the
# subjects that have recorded values.
$ws.Range("R3").Value = 7
$ws.Range("S3").Value = 8

$ws.Range("R5").Value = 8
$ws.Range("S5").Value = 9

$ws.Range("R6").Value = 8
$ws.Range("S6").Value = 9

$ws.Range("R7").Value = 9
$ws.Range("S7").Value = 10

$ws.Range("R8").Value = 9
$ws.Range("S8").Value = 12

$ws.Range("R9").Value = 12
$ws.Range("S9").Value = 15

$ws.Range("R10").Value = 11
$ws.Range("S10").Value = 16

$ws.Range("R11").Value = 11
$ws.Range("S11").Value = 14

$ws.Range("R12").Value = 12
$ws.Range("S12").Value = 11.5

$ws.Range("R13").Value = 10
$ws.Range("S13").Value = 10

$ws.Range("R14").Value = 9
$ws.Range("S14").Value = 10

$ws.Range("R15").Value = 11
$ws.Range("S15").Value = 15

$ws.Range("R16").Value = 12
$ws.Range("S16").Value = 12

$ws.Range("R17").Value = 14
$ws.Range("S17").Value = 24

$ws.Range("R18").Value = 12
$ws.Range("S18").Value = 15

$ws.Range("R19").Value = 10
$ws.Range("S19").Value = 13

$ws.Range("R20").Value = 10
$ws.Range("S20").Value = 12

$ws.Range("R21").Value = 8
$ws.Range("S21").Value = 8

$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 10

$ws.Range("R23").Value = 9
$ws.Range("S23").Value = 12

$ws.Range("R24").Value = 10
$ws.Range("S24").Value = 10

# Leave the cursor where the author last left it.
$ws.Range("U21").Select()
